$wb = $excel.ActiveWorkbook

$wsTable  = $wb.Worksheets.Item("Table")
$wsInput  = $wb.Worksheets.Item("Input")
$wsOutput = $wb.Worksheets.Item("Output")

# The per-car "1: Nissan / 2:Chev / 3:Tesla" helper labels in the Input
# sheet (F2:F4) were stale leftovers - clear them out.
$wsInput.Range("F2:F4").ClearContents()

# Move (cut/paste) the little charging-point lookup table that lived
# below the customer data on Input (A20:F25) onto the Output sheet,
# right under the results table, starting at A12.
$wsInput.Range("A20:F25").Cut($wsOutput.Range("A12"))

# Remove the now-empty rows from Input so the sheet shrinks back down.
$wsInput.Range("A20:F25").EntireRow.Delete()

# Leave the same kind of selection state an end user would have after
# doing the cut (Input) / paste (Output) by hand.
$wsInput.Select()
$wsInput.Range("A20:F25").Select()

$wsOutput.Select()
$wsOutput.Range("A12:F17").Select()

$wsTable.Range("C3").Select()

# The Output sheet is where the user ended up / is now the active tab.
$wsOutput.Activate()
